$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.676.42"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.287.36"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "103.66"
$ws.Range("E5").Value = "  +6.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.45"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -3.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.86"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0932"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.96"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.55"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.855"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("D16").Value = "2.288.47"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").Value = "43.657.94"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.24"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("E21").Value = "  +9.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "233.17"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("E23").Value = "  +13.57%  "
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.22"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.48"
$ws.Range("E27").Value = "  +6.37%  "
$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.45"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "177.25"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.77"
$ws.Range("E31").Value = "  -2.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0899"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.89"
$ws.Range("E34").Value = "  +11.53%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("E38").Value = "  +6.30%  "
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.28"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.49"
$ws.Range("E43").Value = "  +5.14%  "
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("E45").Value = "  -4.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.101"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("E47").Value = "  +2.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.96"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.444"
$ws.Range("E49").Value = "  +7.43%  "
$ws.Range("E50").Value = "  +10.74%  "
$ws.Range("D51").Value = "2.511.40"
$ws.Range("E51").Value = "  -1.11%  "
